$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestSheet")

# The database values can now be edited even when the workbook isn't the
# active/open one - update the second row and append a new row of data.
$ws.Range("A2").Value = "Jacob"
$ws.Range("A4").Value = "Eleanor"
$ws.Range("B4").Value = "Rigsby"

# Keep the named range in sync with the new extent of the table.
$n = $wb.Names.Item("TestName")
$n.RefersTo = '=TestSheet!$A$1:$B$4'

# Make sure the workbook is fully recalculated the next time it is opened,
# since edits are now possible while it is closed.
$wb.ForceFullCalculation = $true
